$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows (2-5) are re-sorted by ascending "Fecha" (column D), which
# shuffles the Volumen (M), Precio mínimo (N), Precio máximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) values along with it.

# New values (row 2..5), ordered ascending by date
$data = @(
    @{ D = 44257; M = 100; N = 14000; O = 15000; P = 14500; S = 806 },
    @{ D = 44252; M = 120; N = 13000; O = 14000; P = 13500; S = 750 },
    @{ D = 44253; M = 160; N = 14000; O = 15000; P = 14500; S = 806 },
    @{ D = 44250; M = 200; N = 14000; O = 15000; P = 14500; S = 806 }
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}

$wb.Save()
